# Add 2022-Q4 data
# --------------------------------------------------------------------------
# 1) Create the new "2022-Q4" worksheet by duplicating the existing
#    "2022-Q3" worksheet (so that column widths, header styles, number
#    formats, page margins, etc. all match exactly) and placing it right
#    before "2022-Q3" - matching the position Q4 should occupy.
# --------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook

$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# The duplicated sheet has 3 data rows (copied from Q3); we only need one,
# so remove the extra two rows, leaving header (row 1) + 1 data row (row 2).
$q4.Rows(3).Delete()
$q4.Rows(3).Delete()

# Fill in the Q4 fund-holding data on row 2.
$q4.Range("A2").Value = 0

$q4.Range("B2").NumberFormat = "@"
$q4.Range("B2").Value = "100055"
$q4.Range("B2").Style = "Normal"

$q4.Range("C2").Value = "富国全球科技互联网股票（QDII）"

$q4.Range("D2").NumberFormat = "@"
$q4.Range("D2").Value = "3.86"
$q4.Range("D2").Style = "Normal"

$q4.Range("E2").NumberFormat = "@"
$q4.Range("E2").Value = "94.32"
$q4.Range("E2").Style = "Normal"

$q4.Range("F2").NumberFormat = "@"
$q4.Range("F2").Value = "3.98"
$q4.Range("F2").Style = "Normal"

$q4.Range("G2").NumberFormat = "@"
$q4.Range("G2").Value = "0.1536"
$q4.Range("G2").Style = "Normal"

$q4.Range("H2").Value = 8

# --------------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert a new "2022-Q4" row right
#    after the header, pushing the existing 2022-Q3 / 2022-Q2 rows down by
#    one, and renumbering the running index in column A.
# --------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Shift existing row 3 (2022-Q2) down to row 4, and row 2 (2022-Q3) down to
# row 3, copying full rows so styles move together with the values.
$total.Range("A3:D3").Copy($total.Range("A4:D4"))
$total.Range("A2:D2").Copy($total.Range("A3:D3"))

# New row 2: 2022-Q4 summary data.
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.15

# Fix up the running index in column A for the rows that moved down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2

# --------------------------------------------------------------------------
# 3) Restore "2022-Q2" as the active/selected sheet (it was the active tab
#    before the edit; adding/copying sheets above shifted the active tab to
#    the newly inserted sheet).
# --------------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Activate()
